$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 98946
$ws.Cells.Item(2, 2).Value = "Raquel da Cruz"
$ws.Cells.Item(2, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 45088
$ws.Cells.Item(2, 7).Value = 12437.3

# Row 3
$ws.Cells.Item(3, 1).Value = 67893
$ws.Cells.Item(3, 2).Value = "Benício Pinto"
$ws.Cells.Item(3, 3).Value = "Engenharia"
$ws.Cells.Item(3, 5).Value = 6
$ws.Cells.Item(3, 6).Value = 45093
$ws.Cells.Item(3, 7).Value = 12398.61

# Row 4
$ws.Cells.Item(4, 1).Value = 40722
$ws.Cells.Item(4, 2).Value = "Nathan Gonçalves"
$ws.Cells.Item(4, 3).Value = "Operações"
$ws.Cells.Item(4, 4).Value = "Outros"
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 45096
$ws.Cells.Item(4, 7).Value = 2997.42

# Row 5
$ws.Cells.Item(5, 1).Value = 93390
$ws.Cells.Item(5, 2).Value = "Bianca Pinto"
$ws.Cells.Item(5, 3).Value = "Marketing"
$ws.Cells.Item(5, 4).Value = "Outros"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 45081
$ws.Cells.Item(5, 7).Value = 7687.89

# Row 6
$ws.Cells.Item(6, 1).Value = 75769
$ws.Cells.Item(6, 2).Value = "Dra. Elisa Campos"
$ws.Cells.Item(6, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(6, 4).Value = "Consulta médica"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 45101
$ws.Cells.Item(6, 7).Value = 11026.77

# Row 7
$ws.Cells.Item(7, 1).Value = 10787
$ws.Cells.Item(7, 2).Value = "Luiz Gustavo Gonçalves"
$ws.Cells.Item(7, 3).Value = "Marketing"
$ws.Cells.Item(7, 4).Value = "Viagem de negócios"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 45098
$ws.Cells.Item(7, 7).Value = 7837.71

# Row 8
$ws.Cells.Item(8, 1).Value = 48662
$ws.Cells.Item(8, 2).Value = "Júlia da Luz"
$ws.Cells.Item(8, 3).Value = "Marketing"
$ws.Cells.Item(8, 4).Value = "Viagem de negócios"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 45083
$ws.Cells.Item(8, 7).Value = 3847.54

# Row 9
$ws.Cells.Item(9, 1).Value = 85395
$ws.Cells.Item(9, 2).Value = "Ana Júlia Rezende"
$ws.Cells.Item(9, 3).Value = "Marketing"
$ws.Cells.Item(9, 4).Value = "Outros"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 45103
$ws.Cells.Item(9, 7).Value = 12272.87

# Row 10
$ws.Cells.Item(10, 1).Value = 24470
$ws.Cells.Item(10, 2).Value = "Elisa Araújo"
$ws.Cells.Item(10, 3).Value = "Vendas"
$ws.Cells.Item(10, 4).Value = "Problemas pessoais"
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 45097
$ws.Cells.Item(10, 7).Value = 4022.48

# Row 11
$ws.Cells.Item(11, 1).Value = 86018
$ws.Cells.Item(11, 2).Value = "Anthony da Rocha"
$ws.Cells.Item(11, 3).Value = "P&D"
$ws.Cells.Item(11, 4).Value = "Consulta médica"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 45094
$ws.Cells.Item(11, 7).Value = 2538.81
